$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Station": drop the (redundant, border-only) explicit cell
# style overrides on V1/V2 so they fall back to the row/column
# default formatting.
# -----------------------------------------------------------------
$wsStation = $wb.Worksheets.Item("Station")
$wsStation.Range("V1").Borders.LineStyle = -4142
$wsStation.Range("V2").Borders.LineStyle = -4142

# -----------------------------------------------------------------
# Sheet "Samples onboard": add a "Weight_g" column, rename the
# length/weight-total headers to carry explicit units, and retire
# the old "Total weight" tallies (the column keeps its header but
# the data is no longer populated; "Number" stays the empty column
# it always was).
# -----------------------------------------------------------------
$wsSamples = $wb.Worksheets.Item("Samples onboard")

# Insert a new blank column before D ("Sex" and everything to its
# right shifts one column to the right: D->E, E->F, F->G, G->H).
$wsSamples.Range("D1:D5").EntireColumn.Insert()

# After the insert, G (old "Total weight") carries the "last column"
# look (border on the right) while H (old "Number") carries the
# "continues" look. Swap the formatting back so the rightmost header
# ("Total weight_g") looks like an interior column and "Number"
# keeps the original trailing-edge look, using a scratch column as a
# buffer.
$wsSamples.Range("G1:G5").Copy()
$wsSamples.Range("Z1:Z5").PasteSpecial(-4122)
$wsSamples.Range("H1:H5").Copy()
$wsSamples.Range("G1:G5").PasteSpecial(-4122)
$wsSamples.Range("Z1:Z5").Copy()
$wsSamples.Range("H1:H5").PasteSpecial(-4122)
$wsSamples.Range("Z1:Z5").Clear()

# Rename headers.
$wsSamples.Range("C1").Value = "Lenght_mm"
$wsSamples.Range("D1").Value = "Weight_g"
$wsSamples.Range("G1").Value = "Number"
$wsSamples.Range("H1").Value = "Total weight_g"

# The old "Total weight" tallies (now sitting in column G after the
# insert/shuffle above) are no longer recorded.
$wsSamples.Range("G2:G5").ClearContents()
